$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each of these cells holds plain text in the source workbook (t="inlineStr"),
# e.g. "304.45" or "0.91%". If we just set .Value, Excel auto-detects these
# as Number / Percentage and converts them (and stamps a new number format),
# which the source diff does not do. So for every touched cell we: mark it
# Text ("@") before writing, write the literal new text, then reset the style
# back to Normal/General so no stray number-format/style sticks around.
$cells = @(
    @{ Addr = "D2"; Value = "304.45" },
    @{ Addr = "E2"; Value = "0.91%" },
    @{ Addr = "D3"; Value = "35.84" },
    @{ Addr = "E3"; Value = "1.36%" },
    @{ Addr = "D4"; Value = "5.068" },
    @{ Addr = "E4"; Value = "-0.12%" },
    @{ Addr = "D5"; Value = "0.08046" },
    @{ Addr = "E5"; Value = "1.46%" },
    @{ Addr = "D6"; Value = "1.924" },
    @{ Addr = "E6"; Value = "1.97%" },
    @{ Addr = "D7"; Value = "4.150" },
    @{ Addr = "E7"; Value = "2.51%" },
    @{ Addr = "D8"; Value = "7.848" },
    @{ Addr = "E8"; Value = "0.89%" },
    @{ Addr = "D9"; Value = "0.9305" },
    @{ Addr = "E9"; Value = "0.16%" },
    @{ Addr = "D10"; Value = "0.1279" },
    @{ Addr = "E10"; Value = "-7.16%" },
    @{ Addr = "D11"; Value = "0.1920" },
    @{ Addr = "E11"; Value = "1.10%" },
    @{ Addr = "D12"; Value = "0.09173" },
    @{ Addr = "E12"; Value = "1.02%" },
    @{ Addr = "D13"; Value = "0.03485" },
    @{ Addr = "E13"; Value = "1.28%" },
    @{ Addr = "D14"; Value = "0.09880" },
    @{ Addr = "D15"; Value = "0.001417" },
    @{ Addr = "E15"; Value = "0.56%" },
    @{ Addr = "D16"; Value = "0.006653" },
    @{ Addr = "E16"; Value = "13.35%" },
    @{ Addr = "D17"; Value = "3.614" },
    @{ Addr = "E17"; Value = "2.34%" },
    @{ Addr = "D18"; Value = "3.050" },
    @{ Addr = "E18"; Value = "2.29%" },
    @{ Addr = "D19"; Value = "0.3421" },
    @{ Addr = "E19"; Value = "-0.13%" },
    @{ Addr = "D20"; Value = "0.1337" },
    @{ Addr = "E20"; Value = "2.51%" },
    @{ Addr = "D21"; Value = "5.180" },
    @{ Addr = "E21"; Value = "2.95%" },
    @{ Addr = "D22"; Value = "0.2534" },
    @{ Addr = "E22"; Value = "5.51%" },
    @{ Addr = "D23"; Value = "0.04405" },
    @{ Addr = "E23"; Value = "-2.07%" },
    @{ Addr = "D24"; Value = "0.001237" },
    @{ Addr = "E24"; Value = "1.73%" },
    @{ Addr = "D25"; Value = "0.004719" },
    @{ Addr = "E25"; Value = "-0.81%" },
    @{ Addr = "E26"; Value = "5.78%" },
    @{ Addr = "D27"; Value = "0.0003133" },
    @{ Addr = "E27"; Value = "4.34%" },
    @{ Addr = "D39"; Value = "0.01994" },
    @{ Addr = "E39"; Value = "7.63%" },
    @{ Addr = "D40"; Value = "0.05158" },
    @{ Addr = "E40"; Value = "8.56%" },
    @{ Addr = "D41"; Value = "0.007604" },
    @{ Addr = "E41"; Value = "3.68%" },
    @{ Addr = "D42"; Value = "0.01006" },
    @{ Addr = "E42"; Value = "3.42%" },
    @{ Addr = "D43"; Value = "0.1363" },
    @{ Addr = "E43"; Value = "2.87%" },
    @{ Addr = "D44"; Value = "0.002103" },
    @{ Addr = "E44"; Value = "-0.39%" },
    @{ Addr = "D45"; Value = "0.01073" },
    @{ Addr = "E45"; Value = "-2.20%" },
    @{ Addr = "D46"; Value = "0.00006394" },
    @{ Addr = "E46"; Value = "2.54%" },
    @{ Addr = "E47"; Value = "0.04%" },
    @{ Addr = "D48"; Value = "63.57" },
    @{ Addr = "E48"; Value = "-1.69%" },
    @{ Addr = "D49"; Value = "0.001603" },
    @{ Addr = "E49"; Value = "-3.43%" },
    @{ Addr = "E50"; Value = "0.04%" },
    @{ Addr = "E51"; Value = "0.04%" }
)

foreach ($c in $cells) {
    $rng = $ws.Range($c.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $c.Value
    $rng.Style = "Normal"
}
